# "add tabel format baru" - renumber the second table's heading on sheet
# "Bab 4" from "Tabel 4.2.3" to "Tabel 4.2.5" (making room for the two new
# tables, 4.2.6 and 4.2.7, that were introduced alongside this one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bab 4")

$ws.Range("H1").Value = "Tabel 4.2.5"
